# The "reviews_count" column (column E) is empty across the whole sheet and
# is being removed, shifting the remaining columns (reviews_average,
# latitude, longitude, is_permanently_closed, gmaps_link,
# latest_review_date) one position to the left (F:K -> E:J).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(5).Delete()
